$d = $word.ActiveDocument

# --- 1. Title: merge "Circle " + "Broader View" + ": System Objects" runs
#        into a single run reading "Circle Language | Broader View" ---
$titleRng = $d.Content
$titleRng.Find.Execute("Circle Broader View: System Objects") | Out-Null
$titleStart = $titleRng.Start
$titleEnd = $titleRng.End

$tail = $d.Range($titleStart + 7, $titleEnd)   # "Broader View: System Objects"
$tail.Delete()

$head = $d.Range($titleStart, $titleStart + 7) # "Circle "
$head.Text = "Circle Language | Broader View"

# --- 2. Remove the _GoBack bookmark (start & end) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Merge the "." + " " + "System objects without any extensions..." runs
#        into a single run ---
$mergeRng = $d.Content
$mergeRng.Find.Execute(". System objects without any extensions need to be implemented right inside the code base. ") | Out-Null
$mergeStart = $mergeRng.Start
$mergeEnd = $mergeRng.End

$restRange = $d.Range($mergeStart + 1, $mergeEnd)
$restRange.Delete()

$dotRange = $d.Range($mergeStart, $mergeStart + 1)
$dotRange.Text = ". System objects without any extensions need to be implemented right inside the code base. "
